$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.225111487315985,0.2006750712341443,0.1579374492244341,0,1.533339266402763,0.002475189058415209,0,0.9362062281404526,0.1940202047970914,0,0.368258523396733,0.3174805604416733,0,3.801774937353116),
    @(1.121270145119638,0.1846934509052858,0.1577935054873976,0,1.544679415121983,0.002477993517129873,0,0.9489868317012231,0.1954383489362055,0,0.3647141112970473,0.3005162760387918,0,3.834246426473754),
    @(1.057536925209604,0.1748252893644349,0.1577437270608293,0,1.552527625682515,0.0024798082109431,0,0.9573873776198987,0.196361499002272,0,0.3626574857707965,0.2901542426530241,0,3.856571794901555),
    @(1.031573406753751,0.1707902509213852,0.1577331865717788,0,1.555948446384313,0.002480571103935357,0,0.9609497030683656,0.1967508923840118,0,0.3618496127635993,0.2859455529281831,0,3.866269687236027),
    @(1.02726273806104,0.1701194166899143,0.1577320259859505,0,1.556529915987745,0.002480699196641567,0,0.9615496197180384,0.1968163489747816,0,0.3617172950190977,0.2852475535097696,0,3.867916252225427),
    @(1.057186736171104,0.1747709264993773,0.1577435454025675,0,1.552572858735132,0.002479818404848191,0,0.9574348575380363,0.1963666970097453,0,0.3626464679860959,0.2900974259951496,0,3.856700154588523),
    @(1.189302686181975,0.1951762440750713,0.1578798273994053,0,1.537065597073628,0.00247613682819277,0,0.9404981201725846,0.1944983194723378,0,0.367011651011282,0.311620229938903,0,3.81247535253226),
    @(1.44852027938424,0.2347427671507489,0.1584518437378364,0,1.513681559714712,0.002469650000985923,0,0.9116776412100123,0.1912491025817076,0,0.3765162238181432,0.3542440023846893,0,3.744711058976833),
    @(1.638982757266547,0.263529984633891,0.1590559976427741,0,1.500786903740817,0.002465326398685249,0,0.8931846853836589,0.1891131047604606,0,0.3840693957238841,0.3858019164046809,0,3.706505970381727),
    @(1.725619096886589,0.2765630148292644,0.1593703832989917,0,1.49585200067115,0.002463454573369717,0,0.8853548047959166,0.1881955944174241,0,0.3876282208110382,0.4002086380611942,0,3.69164528115229),
    @(1.758423638306795,0.2814891047784727,0.1594950855071815,0,1.494117189634935,0.002462759352079846,0,0.8824737010821089,0.1878559213664719,0,0.38899340661861,0.4056711332700971,0,3.686380542409353),
    @(1.751358746192011,0.2804285980461714,0.1594679778971724,0,1.494484854704751,0.002462908476613471,0,0.8830904648907563,0.1879287308704782,0,0.3886986115690121,0.4044943819171323,0,3.687498260546135),
    @(1.728318012391526,0.2769684737929197,0.1593805295748041,0,1.495706592224586,0.002463397105020535,0,0.8851160924682198,0.1881674937425064,0,0.3877401847619382,0.4006579027887085,0,3.691204876661459),
    @(1.714204482109835,0.2748478373354146,0.1593276998231445,0,1.496472385135874,0.002463698172952562,0,0.8863677778417376,0.188314753933998,0,0.3871554005349367,0.3983088480298989,0,3.693522531319985),
    @(1.633320585505203,0.2626769671554143,0.1590362445849109,0,1.501128149449528,0.002465450633201921,0,0.893708121034237,0.1891741543911536,0,0.3838392790856062,0.3848613983983,0,3.707527878830774),
    @(1.583698044115067,0.2551943598049604,0.1588675506902959,0,1.504222790374939,0.002466549998314329,0,0.8983604995232817,0.1897152270262019,0,0.3818363182504498,0.376624615001333,0,3.716765204821797),
    @(1.55515601887879,0.250884702155389,0.1587742475126532,0,1.506090366900402,0.002467191270070211,0,0.9010912601191592,0.1900315372588857,0,0.3806958390250088,0.3718918427439775,0,3.722315343445302),
    @(1.545492159994524,0.2494245276492961,0.1587432975618484,0,1.506737743311618,0.00246740993189809,0,0.9020252632379595,0.1901395109608082,0,0.380311683621315,0.3702902432684496,0,3.724235228697268),
    @(1.588980508917814,0.2559915044902255,0.1588851231764608,0,1.50388429167883,0.002466432043725531,0,0.8978595697710681,0.1896571012861585,0,0.3820483402497672,0.3775009394913837,0,3.71575733623385),
    @(1.73508572438999,0.2779850488760758,0.1594060621852904,0,1.495344103001848,0.002463253214602431,0,0.8845188385381384,0.1880971526371713,0,0.3880212230563131,0.4017845828450106,0,3.690106306984859),
    @(1.830556969550457,0.2923051495134814,0.1597794485578419,0,1.490543268394418,0.002461254898300994,0,0.8762889993392982,0.1871229056033474,0,0.3920270006377109,0.4176958934873838,0,3.675456071417756),
    @(1.779604344503582,0.2846672642363046,0.1595771648128945,0,1.493034112527454,0.002462314208258111,0,0.8806366282372302,0.1876387442026921,0,0.3898797373867922,0.4092001268736354,0,3.683081579995957),
    @(1.586592348068336,0.2556311396776039,0.1588771671867804,0,1.50403705150044,0.002466485342308469,0,0.8980858655385617,0.1896833636083528,0,0.3819524506762946,0.3771047448658393,0,3.71621224769504),
    @(1.378387420023955,0.2240879053355229,0.1582646678893198,0,1.519255148204074,0.002471326874499742,0,0.9190037219731764,0.1920838743794366,0,0.3738444557761653,0.342669720181668,0,3.761010911946329)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 2).Value = $rowVals[$j]
    }
}